# Generate Report for Handback
# Row 7 (the d7c2841e-8762-4c7c-95aa-1293904d31ff handback) has finished
# processing on both the zh-cn and de-de sheets: fill in the "Latest
# Target File" / "Latest Handback File" / "Latest Handback DateTime" /
# "Error Detail" columns, matching the pattern already used for the
# earlier rows.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3dfc925457a3a87ebb1a675e0edb500507e27b0/e2e/d7c2841e-8762-4c7c-95aa-1293904d31ff.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6991b0a10c1506646c8dd7914af3a64b85cf43a4/e2e/d7c2841e-8762-4c7c-95aa-1293904d31ff.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3dfc925457a3a87ebb1a675e0edb500507e27b0/e2e/d7c2841e-8762-4c7c-95aa-1293904d31ff.md."

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("J7").Value = "d7c2841e-8762-4c7c-95aa-1293904d31ff.0fd5ce66175484c4b4c9accf052f126bdbcbff95.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-26 18:57:41"
$wsZh.Range("P7").Value = $errorDetail
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, [System.Type]::Missing, "d7c2841e-8762-4c7c-95aa-1293904d31ff.md", "d7c2841e-8762-4c7c-95aa-1293904d31ff.md")

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("J7").Value = "d7c2841e-8762-4c7c-95aa-1293904d31ff.0fd5ce66175484c4b4c9accf052f126bdbcbff95.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-26 18:57:48"
$wsDe.Range("P7").Value = $errorDetail
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, [System.Type]::Missing, "d7c2841e-8762-4c7c-95aa-1293904d31ff.md", "d7c2841e-8762-4c7c-95aa-1293904d31ff.md")
